$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2549.5
$ws.Range("J17").Value = 2549.5
$ws.Range("L17").Value = 7648.5
$ws.Range("N17").Value = -7984.5
$ws.Range("H33").Value = 296.2
$ws.Range("I33").Value = 308
$ws.Range("J33").Value = 249
$ws.Range("K33").Value = 308
$ws.Range("L33").Value = 249
$ws.Range("M33").Value = -79
$ws.Range("N33").Value = -707
$ws.Range("H51").Value = 1500
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 1500
$ws.Range("M51").Value = -1016
$ws.Range("H64").Value = 3861
$ws.Range("I64").Value = 3861
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3861
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -3613
$ws.Range("H67").Value = 3861
$ws.Range("I67").Value = 3861
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3861
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -3003
$ws.Range("H86").Value = 12156
$ws.Range("I86").Value = 6999.6665
$ws.Range("J86").Value = 15249.8
$ws.Range("K86").Value = 6999.6665
$ws.Range("L86").Value = 15249.8
$ws.Range("M86").Value = -5876.6665
$ws.Range("N86").Value = -17495.8
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = ""
$ws.Range("N87").Value = 0
$ws.Range("H89").Value = 12156
$ws.Range("I89").Value = 6999.6665
$ws.Range("J89").Value = 15249.8
$ws.Range("K89").Value = 34998.3325
$ws.Range("L89").Value = 76249
$ws.Range("M89").Value = -29382.3325
$ws.Range("N89").Value = -87481
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = ""
$ws.Range("N90").Value = 0
$ws.Range("H98").Value = 690.4706
$ws.Range("I98").Value = 438.13333
$ws.Range("K98").Value = 438.13333
$ws.Range("M98").Value = 1059.86667
$ws.Range("H116").Value = 5750
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 5750
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = ""
$ws.Range("M116").Value = 5750
$ws.Range("N116").Value = -12634
$ws.Range("H122").Value = 690.4706
$ws.Range("I122").Value = 438.13333
$ws.Range("K122").Value = 1314.39999
$ws.Range("M122").Value = 1135.60001
$ws.Range("H125").Value = 4665.5
$ws.Range("J125").Value = 8474.5
$ws.Range("L125").Value = 76270.5
$ws.Range("N125").Value = -81190.5
$ws.Range("H135").Value = 1083
$ws.Range("I135").Value = 776.5
$ws.Range("K135").Value = 6988.5
$ws.Range("M135").Value = -4453.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3420
$ws.Range("I2").Value = 1130.25
$ws.Range("K2").Value = 1130.25
$ws.Range("M2").Value = -1017.25
$ws.Range("H45").Value = 1822.6666
$ws.Range("I45").Value = 1609
$ws.Range("K45").Value = 1609
$ws.Range("M45").Value = -1232
$ws.Range("H92").Value = 69889
$ws.Range("J92").Value = 69889
$ws.Range("L92").Value = 69889
$ws.Range("N92").Value = -74881
$ws.Range("H94").Value = 90330
$ws.Range("J94").Value = 90330
$ws.Range("L94").Value = 90330
$ws.Range("N94").Value = -92132
$ws.Range("H116").Value = 3420
$ws.Range("I116").Value = 1130.25
$ws.Range("K116").Value = 1130.25
$ws.Range("M116").Value = 1163.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3420
$ws.Range("I3").Value = 1130.25
$ws.Range("K3").Value = 1130.25
$ws.Range("M3").Value = -1016.25
$ws.Range("H95").Value = 19049
$ws.Range("J95").Value = 19049
$ws.Range("L95").Value = 19049
$ws.Range("N95").Value = -24541
$ws.Range("H99").Value = 1767.1904
$ws.Range("I99").Value = 1945.1111
$ws.Range("J99").Value = 699.6667
$ws.Range("K99").Value = 1945.1111
$ws.Range("L99").Value = 699.6667
$ws.Range("M99").Value = -447.1111000000001
$ws.Range("N99").Value = -3695.6667
$ws.Range("H105").Value = 1246.4
$ws.Range("I105").Value = 699
$ws.Range("K105").Value = 699
$ws.Range("M105").Value = 1048
$ws.Range("H134").Value = 4698.6787
$ws.Range("I134").Value = 4698.6787
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14096.0361
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = -11561.0361

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 528
$ws.Range("I16").Value = 671.6
$ws.Range("K16").Value = 671.6
$ws.Range("M16").Value = -384.6
$ws.Range("H58").Value = 2184.8333
$ws.Range("I58").Value = 2184.8333
$ws.Range("K58").Value = 2184.8333
$ws.Range("M58").Value = -1981.8333
$ws.Range("H75").Value = 24800
$ws.Range("J75").Value = 24800
$ws.Range("L75").Value = 24800
$ws.Range("N75").Value = -26796
$ws.Range("H78").Value = 24800
$ws.Range("J78").Value = 24800
$ws.Range("L78").Value = 74400
$ws.Range("N78").Value = -84384
$ws.Range("H88").Value = 26661.857
$ws.Range("J88").Value = 19438.834
$ws.Range("L88").Value = 19438.834
$ws.Range("N88").Value = -20250.834
$ws.Range("H91").Value = 26661.857
$ws.Range("J91").Value = 19438.834
$ws.Range("L91").Value = 19438.834
$ws.Range("N91").Value = -22246.834
$ws.Range("H113").Value = 528
$ws.Range("I113").Value = 671.6
$ws.Range("K113").Value = 671.6
$ws.Range("M113").Value = 1498.4
$ws.Range("H132").Value = 2522.1428
$ws.Range("J132").Value = 2935
$ws.Range("L132").Value = 8805
$ws.Range("N132").Value = -13865
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139
$ws.Range("H136").Value = 2184.8333
$ws.Range("I136").Value = 2184.8333
$ws.Range("K136").Value = 6554.499899999999
$ws.Range("M136").Value = -4004.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5327.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5327.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = ""
$ws.Range("M80").Value = 5327.5
$ws.Range("N80").Value = -7323.5
$ws.Range("H83").Value = 5327.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5327.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = ""
$ws.Range("M83").Value = 26637.5
$ws.Range("N83").Value = -36621.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4293.875
$ws.Range("I7").Value = 4193
$ws.Range("K7").Value = 4193
$ws.Range("M7").Value = -4081
$ws.Range("H68").Value = 1326
$ws.Range("I68").Value = 1326
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1326
$ws.Range("L68").Value = ""
$ws.Range("M68").Value = -577
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 1326
$ws.Range("I71").Value = 1326
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6630
$ws.Range("L71").Value = ""
$ws.Range("M71").Value = -2886
$ws.Range("N71").Value = 0
$ws.Range("H82").Value = 1874.875
$ws.Range("I82").Value = 1299.8
$ws.Range("J82").Value = 2833.3333
$ws.Range("K82").Value = 1299.8
$ws.Range("L82").Value = 2833.3333
$ws.Range("M82").Value = -938.8
$ws.Range("N82").Value = -3555.3333
$ws.Range("H85").Value = 1874.875
$ws.Range("I85").Value = 1299.8
$ws.Range("J85").Value = 2833.3333
$ws.Range("K85").Value = 1299.8
$ws.Range("L85").Value = 2833.3333
$ws.Range("M85").Value = -51.79999999999995
$ws.Range("N85").Value = -5329.3333
$ws.Range("H93").Value = 1151
$ws.Range("I93").Value = 604
$ws.Range("J93").Value = 1333.3334
$ws.Range("K93").Value = 604
$ws.Range("L93").Value = 1333.3334
$ws.Range("M93").Value = 644
$ws.Range("N93").Value = -3829.3334
$ws.Range("H94").Value = 56799.8
$ws.Range("J94").Value = 56799.8
$ws.Range("L94").Value = 56799.8
$ws.Range("N94").Value = -58151.8
$ws.Range("H126").Value = 4293.875
$ws.Range("I126").Value = 4193
$ws.Range("K126").Value = 12579
$ws.Range("M126").Value = -10109

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8333.333000000001
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5336
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = ""
$ws.Range("N95").Value = 0
$ws.Range("H97").Value = 42499
$ws.Range("J97").Value = 42499
$ws.Range("L97").Value = 42499
$ws.Range("N97").Value = -44481
$ws.Range("H126").Value = 1032.4
$ws.Range("I126").Value = 891.8570999999999
$ws.Range("K126").Value = 2675.5713
$ws.Range("M126").Value = -205.5712999999996
$ws.Range("H132").Value = 2241.1875
$ws.Range("I132").Value = 2143.8462
$ws.Range("J132").Value = 2663
$ws.Range("K132").Value = 6431.5386
$ws.Range("L132").Value = 7989
$ws.Range("M132").Value = -3901.5386
$ws.Range("N132").Value = -13049
$ws.Range("H136").Value = 3987.2307
$ws.Range("I136").Value = 2894.0908
$ws.Range("J136").Value = 9999.5
$ws.Range("K136").Value = 8682.2724
$ws.Range("L136").Value = 29998.5
$ws.Range("M136").Value = -6132.2724
$ws.Range("N136").Value = -35098.5
